# Revert "Powerpoint writer: consolidate text run nodes."
#
# Split previously-merged text runs back into separate <a:r> nodes so that
# trailing/leading space characters live in their own run, matching the
# original (pre-consolidation) run layout.
#
# Mechanism: re-assigning TextRange.Characters(start, length).Text forces
# the writer to break the enclosing run at that character boundary without
# changing the visible text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title 1: "A " + "slide"  ->  "A" + " " + "slide" ---
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Characters(1, 1).Text = "A"

# --- TextBox 3: "Followed " + "by " + "a " + "picture"
#              ->  "Followed" + " " + "by" + " " + "a" + " " + "picture" ---
$caption = $s.Shapes.Item(4)
$captionRange = $caption.TextFrame.TextRange
$captionRange.Characters(9, 1).Text = " "
$captionRange.Characters(12, 1).Text = " "
$captionRange.Characters(14, 1).Text = " "
